# Updated symbol list: refreshed Price/Volume(1h) figures, plus a coin-order
# reshuffle in rows 15-19 (One/TigerCash/LEO/GateToken/BTSEToken), per the
# "Updated symbol list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All the data in this sheet is stored as TEXT (not numbers), even the price
# and percentage columns. Plain "$range.Value = ..." lets Excel auto-detect
# numeric-looking strings and silently convert them to real numbers, so every
# write below forces Text format first (and clears it back to Normal after)
# to keep the cell a string, matching the original formatting.
function Set-Text($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
$de = $ws.Range("D2:E2")
$de.NumberFormat = "@"
$ws.Range("D2").Value = "258.39"
$ws.Range("E2").Value = "0.27%"
$de.Style = "Normal"

# Row 3
$de = $ws.Range("D3:E3")
$de.NumberFormat = "@"
$ws.Range("D3").Value = "27.04"
$ws.Range("E3").Value = "-0.34%"
$de.Style = "Normal"

# Row 4
$de = $ws.Range("D4:E4")
$de.NumberFormat = "@"
$ws.Range("D4").Value = "4.656"
$ws.Range("E4").Value = "-5.26%"
$de.Style = "Normal"

# Row 5
$de = $ws.Range("D5:E5")
$de.NumberFormat = "@"
$ws.Range("D5").Value = "0.05908"
$ws.Range("E5").Value = "-0.70%"
$de.Style = "Normal"

# Row 6
Set-Text $ws.Range("E6") "-0.76%"

# Row 7
$de = $ws.Range("D7:E7")
$de.NumberFormat = "@"
$ws.Range("D7").Value = "0.8547"
$ws.Range("E7").Value = "-1.60%"
$de.Style = "Normal"

# Row 8
$de = $ws.Range("D8:E8")
$de.NumberFormat = "@"
$ws.Range("D8").Value = "0.9508"
$ws.Range("E8").Value = "-0.53%"
$de.Style = "Normal"

# Row 9
$de = $ws.Range("D9:E9")
$de.NumberFormat = "@"
$ws.Range("D9").Value = "0.1404"
$ws.Range("E9").Value = "-0.59%"
$de.Style = "Normal"

# Row 10
$de = $ws.Range("D10:E10")
$de.NumberFormat = "@"
$ws.Range("D10").Value = "0.05327"
$ws.Range("E10").Value = "50.00%"
$de.Style = "Normal"

# Row 11
$de = $ws.Range("D11:E11")
$de.NumberFormat = "@"
$ws.Range("D11").Value = "0.07092"
$ws.Range("E11").Value = "-1.14%"
$de.Style = "Normal"

# Row 12
Set-Text $ws.Range("E12") "-1.21%"

# Row 13
$de = $ws.Range("D13:E13")
$de.NumberFormat = "@"
$ws.Range("D13").Value = "0.09151"
$ws.Range("E13").Value = "-1.12%"
$de.Style = "Normal"

# Row 14
$de = $ws.Range("D14:E14")
$de.NumberFormat = "@"
$ws.Range("D14").Value = "0.001536"
$ws.Range("E14").Value = "-0.98%"
$de.Style = "Normal"

# Row 15
Set-Text $ws.Range("B15") "One"
Set-Text $ws.Range("C15") "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$de = $ws.Range("D15:E15")
$de.NumberFormat = "@"
$ws.Range("D15").Value = "0.0006029"
$ws.Range("E15").Value = "-0.52%"
$de.Style = "Normal"

# Row 16
Set-Text $ws.Range("B16") "TigerCash"
Set-Text $ws.Range("C16") "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$de = $ws.Range("D16:E16")
$de.NumberFormat = "@"
$ws.Range("D16").Value = "0.006082"
$ws.Range("E16").Value = "1.14%"
$de.Style = "Normal"

# Row 17
Set-Text $ws.Range("B17") "LEO"
Set-Text $ws.Range("C17") "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$de = $ws.Range("D17:E17")
$de.NumberFormat = "@"
$ws.Range("D17").Value = "3.501"
$ws.Range("E17").Value = "0.46%"
$de.Style = "Normal"

# Row 18
Set-Text $ws.Range("B18") "GateToken"
Set-Text $ws.Range("C18") "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$de = $ws.Range("D18:E18")
$de.NumberFormat = "@"
$ws.Range("D18").Value = "3.188"
$ws.Range("E18").Value = "-2.20%"
$de.Style = "Normal"

# Row 19
Set-Text $ws.Range("B19") "BTSEToken"
Set-Text $ws.Range("C19") "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$de = $ws.Range("D19:E19")
$de.NumberFormat = "@"
$ws.Range("D19").Value = "2.204"
$ws.Range("E19").Value = "-0.09%"
$de.Style = "Normal"

# Row 20
$de = $ws.Range("D20:E20")
$de.NumberFormat = "@"
$ws.Range("D20").Value = "0.3055"
$ws.Range("E20").Value = "-2.94%"
$de.Style = "Normal"

# Row 21
Set-Text $ws.Range("E21") "-2.21%"

# Row 22
$de = $ws.Range("D22:E22")
$de.NumberFormat = "@"
$ws.Range("D22").Value = "3.829"
$ws.Range("E22").Value = "8.38%"
$de.Style = "Normal"

# Row 23
Set-Text $ws.Range("E23") "-0.61%"

# Row 24
$de = $ws.Range("D24:E24")
$de.NumberFormat = "@"
$ws.Range("D24").Value = "0.001219"
$ws.Range("E24").Value = "-0.02%"
$de.Style = "Normal"

# Row 25
$de = $ws.Range("D25:E25")
$de.NumberFormat = "@"
$ws.Range("D25").Value = "0.004297"
$ws.Range("E25").Value = "-4.93%"
$de.Style = "Normal"

# Row 27
$de = $ws.Range("D27:E27")
$de.NumberFormat = "@"
$ws.Range("D27").Value = "0.0001937"
$ws.Range("E27").Value = "29.87%"
$de.Style = "Normal"

# Row 40
$de = $ws.Range("D40:E40")
$de.NumberFormat = "@"
$ws.Range("D40").Value = "0.03831"
$ws.Range("E40").Value = "-0.10%"
$de.Style = "Normal"

# Row 41
$de = $ws.Range("D41:E41")
$de.NumberFormat = "@"
$ws.Range("D41").Value = "0.006208"
$ws.Range("E41").Value = "-5.65%"
$de.Style = "Normal"

# Row 42
Set-Text $ws.Range("E42") "0.03%"

# Row 43
$de = $ws.Range("D43:E43")
$de.NumberFormat = "@"
$ws.Range("D43").Value = "0.002339"
$ws.Range("E43").Value = "6.33%"
$de.Style = "Normal"

# Row 44
$de = $ws.Range("D44:E44")
$de.NumberFormat = "@"
$ws.Range("D44").Value = "0.01411"
$ws.Range("E44").Value = "33.95%"
$de.Style = "Normal"

# Row 45
$de = $ws.Range("D45:E45")
$de.NumberFormat = "@"
$ws.Range("D45").Value = "0.00005397"
$ws.Range("E45").Value = "-1.70%"
$de.Style = "Normal"

# Row 46
Set-Text $ws.Range("E46") "-0.02%"

# Row 48
$de = $ws.Range("D48:E48")
$de.NumberFormat = "@"
$ws.Range("D48").Value = "0.2516"
$ws.Range("E48").Value = "11,714.60%"
$de.Style = "Normal"

# Row 49
$de = $ws.Range("D49:E49")
$de.NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("E49").Value = "-0.02%"
$de.Style = "Normal"

# Row 50
$de = $ws.Range("D50:E50")
$de.NumberFormat = "@"
$ws.Range("D50").Value = "0.0002000"
$ws.Range("E50").Value = "-0.02%"
$de.Style = "Normal"
